$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "43.092.96"
$ws.Range("E2").Value = "  +0.13%  "
$ws.Range("D3").Value = "2.304.47"
$ws.Range("E4").Value = "  -0.05%  "
$ws.Range("D5").Value = "'300.28"
$ws.Range("E5").Value = "  -0.11%  "
$ws.Range("D6").Value = "'98.10"
$ws.Range("E6").Value = "  -1.67%  "
$ws.Range("E7").Value = "  +2.90%  "
$ws.Range("E8").Value = "  -0.02%  "
$ws.Range("E9").Value = "  +1.33%  "
$ws.Range("D10").Value = "'36.12"
$ws.Range("E10").Value = "  -0.43%  "
$ws.Range("E11").Value = "  +0.13%  "
$ws.Range("E12").Value = "  +0.56%  "
$ws.Range("D13").Value = "'17.71"
$ws.Range("E13").Value = "  -2.55%  "
$ws.Range("E14").Value = "  -0.49%  "
$ws.Range("D15").Value = "2.662.68"
$ws.Range("E15").Value = "  +0.09%  "
$ws.Range("D16").Value = "2.302.41"
$ws.Range("E16").Value = "  -0.77%  "
$ws.Range("D17").Value = "'0.789"
$ws.Range("E17").Value = "  -1.39%  "
$ws.Range("D18").Value = "42.983.28"
$ws.Range("E18").Value = "  +0.08%  "
$ws.Range("D19").Value = "'12.74"
$ws.Range("E19").Value = "  -0.14%  "
$ws.Range("E20").Value = "  +0.86%  "
$ws.Range("D21").Value = "'6.15"
$ws.Range("E21").Value = "  +0.36%  "
$ws.Range("D22").Value = "'68.48"
$ws.Range("E22").Value = "  +0.87%  "
$ws.Range("D23").Value = "'238.02"
$ws.Range("E24").Value = "  -0.93%  "
$ws.Range("E25").Value = "  -0.40%  "
$ws.Range("E26").Value = "  -0.33%  "
$ws.Range("E27").Value = "  -0.33%  "
$ws.Range("D28").Value = "'25.07"
$ws.Range("E28").Value = "  +0.52%  "
$ws.Range("D29").Value = "'164.11"
$ws.Range("E29").Value = "  -2.10%  "
$ws.Range("D30").Value = "'2.05"
$ws.Range("E30").Value = "  -13.02%  "
$ws.Range("E31").Value = "  +0.18%  "
$ws.Range("D32").Value = "'33.10"
$ws.Range("E32").Value = "  -4.45%  "
$ws.Range("E33").Value = "  -0.03%  "
$ws.Range("E34").Value = "  +1.76%  "
$ws.Range("D35").Value = "'4.82"
$ws.Range("E35").Value = "  +4.14%  "
$ws.Range("D36").Value = "'18.11"
$ws.Range("E36").Value = "  +2.81%  "
$ws.Range("E37").Value = "  +0.29%  "
$ws.Range("E38").Value = "  +1.27%  "
$ws.Range("E39").Value = "  +1.14%  "
$ws.Range("E40").Value = "  -0.30%  "
$ws.Range("E41").Value = "  -0.88%  "
$ws.Range("E42").Value = "  +1.16%  "
$ws.Range("D43").Value = "2.021.83"
$ws.Range("E43").Value = "  +2.11%  "
$ws.Range("D44").Value = "'0.0286"
$ws.Range("E44").Value = "  -1.71%  "
$ws.Range("D45").Value = "'2.23"
$ws.Range("E45").Value = "  -2.89%  "
$ws.Range("D46").Value = "'10.38"
$ws.Range("E46").Value = "  +2.07%  "
$ws.Range("D47").Value = "'17.55"
$ws.Range("E47").Value = "  +0.39%  "
$ws.Range("E48").Value = "  -2.16%  "
$ws.Range("D49").Value = "'54.31"
$ws.Range("E49").Value = "  -1.83%  "
$ws.Range("D50").Value = "2.527.63"
$ws.Range("E50").Value = "  +0.20%  "
$ws.Range("E51").Value = "  -1.06%  "
